$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.944.52"
$ws.Range("E2").Value = "  -1.02%  "
$ws.Range("D3").Value = "1.637.61"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("E8").Value = "  -0.76%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0637"
$ws.Range("D9").Style = "Normal"
$ws.Range("E10").Value = "  -1.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0794"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("D13").Value = "1.863.08"
$ws.Range("E13").Value = "  -0.51%  "
$ws.Range("D14").Value = "1.654.40"
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("E15").Value = "  -0.48%  "
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.86"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "25.902.07"
$ws.Range("E18").Value = "  -1.21%  "
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "192.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.33%  "
$ws.Range("E21").Value = "  -1.94%  "
$ws.Range("E22").Value = "  -1.43%  "
$ws.Range("E23").Value = "  -0.36%  "
$ws.Range("E24").Value = "  +0.72%  "
$ws.Range("B25").Value = "Stellar"
$ws.Range("C25").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.131"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.85%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "144.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.11%  "
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("E28").Value = "  -0.54%  "
$ws.Range("E29").Value = "  -0.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0499"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.41%  "
$ws.Range("E32").Value = "  -2.10%  "
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("E34").Value = "  -3.91%  "
$ws.Range("E35").Value = "  +1.85%  "
$ws.Range("E36").Value = "  -1.37%  "
$ws.Range("D37").Value = "1.132.85"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("E38").Value = "  -1.77%  "
$ws.Range("E39").Value = "  -1.75%  "
$ws.Range("E40").Value = "  -0.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.42"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.15%  "
$ws.Range("E43").Value = "  -0.69%  "
$ws.Range("D44").Value = "1.772.59"
$ws.Range("E44").Value = "  -0.49%  "
$ws.Range("D45").Value = "0.0₆0115"
$ws.Range("E45").Value = "  +3.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.54"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.75%  "
$ws.Range("E47").Value = "  +2.51%  "
$ws.Range("E48").Value = "  -0.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.70"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.415"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0959"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.08%  "
